$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Footer (default) -> footer2.xml: Pearson logo "image2.png" -> "image1.png" ---
$footerDefault = $sec.Footers.Item(1)
$imgFooterDefault = $footerDefault.Range.InlineShapes.Item(1)
$imgFooterDefault.Select() | Out-Null
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# --- Footer (first page) -> footer1.xml: Pearson logo "image2.png" -> "image1.png" ---
$footerFirst = $sec.Footers.Item(2)
$imgFooterFirst = $footerFirst.Range.InlineShapes.Item(1)
$imgFooterFirst.Select() | Out-Null
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# --- Header (first page) -> header1.xml: BTEC logo "image1.jpg" -> "image2.jpg" ---
$headerFirst = $sec.Headers.Item(2)
$imgHeaderFirst = $headerFirst.Range.InlineShapes.Item(1)
$imgHeaderFirst.Select() | Out-Null
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

Write-Output "Renamed inline image shapes successfully."
